$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.626.37'
$ws.Range("E2").Value = '  -1.83%  '

$ws.Range("D3").Value = '1.587.61'
$ws.Range("E3").Value = '  -2.25%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '211.12'
$ws.Range("E5").Value = '  -1.40%  '

$ws.Range("E6").Value = '  -2.43%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '0.247'
$ws.Range("E8").Value = '  -2.11%  '

$ws.Range("E9").Value = '  -1.82%  '

$ws.Range("D10").Value = '19.57'
$ws.Range("E10").Value = '  -3.60%  '

$ws.Range("E11").Value = '  -1.61%  '

$ws.Range("D12").Value = '1.809.72'

$ws.Range("D13").Value = '1.563.51'
$ws.Range("E13").Value = '  -3.79%  '

$ws.Range("D14").Value = '4.03'
$ws.Range("E14").Value = '  -2.84%  '

$ws.Range("E15").Value = '  -4.00%  '

$ws.Range("D16").Value = '64.60'
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").Value = '26.600.96'
$ws.Range("E17").Value = '  -1.85%  '

$ws.Range("E18").Value = '  -2.23%  '

$ws.Range("D19").Value = '208.81'
$ws.Range("E19").Value = '  -3.64%  '

$ws.Range("E20").Value = '  +0.17%  '

$ws.Range("D21").Value = '6.72'
$ws.Range("E21").Value = '  -3.05%  '

$ws.Range("E22").Value = '  -2.67%  '

$ws.Range("D23").Value = '2.35'
$ws.Range("E23").Value = '  -2.41%  '

$ws.Range("D24").Value = '8.87'
$ws.Range("E24").Value = '  -2.11%  '

$ws.Range("D25").Value = '146.76'
$ws.Range("E25").Value = '  -0.90%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -0.75%  '

$ws.Range("E28").Value = '  -2.89%  '

$ws.Range("D29").Value = '15.30'
$ws.Range("E29").Value = '  -2.01%  '

$ws.Range("D30").Value = '0.0508'
$ws.Range("E30").Value = '  +0.20%  '

$ws.Range("E31").Value = '  -1.82%  '

$ws.Range("E32").Value = '  -3.64%  '

$ws.Range("D33").Value = '0.688'
$ws.Range("E33").Value = '  +24.28%  '

$ws.Range("D34").Value = '2.91'
$ws.Range("E34").Value = '  -2.89%  '

$ws.Range("D35").Value = '1.308.72'
$ws.Range("E35").Value = '  -2.60%  '

$ws.Range("E36").Value = '  -0.94%  '

$ws.Range("E37").Value = '  -5.61%  '

$ws.Range("E38").Value = '  -3.22%  '

$ws.Range("E39").Value = '  -3.31%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '5.38'
$ws.Range("E41").Value = '  +2.69%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.789'
$ws.Range("E42").Value = '  -1.90%  '

$ws.Range("E43").Value = '  -2.19%  '

$ws.Range("D44").Value = '62.69'
$ws.Range("E44").Value = '  -4.42%  '

$ws.Range("D45").Value = '1.723.15'

$ws.Range("D46").Value = '89.54'
$ws.Range("E46").Value = '  -1.22%  '

$ws.Range("E47").Value = '  -0.95%  '

$ws.Range("D48").Value = '0.839'
$ws.Range("E48").Value = '  -5.12%  '

$ws.Range("E49").Value = '  -1.61%  '

$ws.Range("D50").Value = '0.0979'
$ws.Range("E50").Value = '  -1.46%  '

$ws.Range("D51").Value = '7.49'
$ws.Range("E51").Value = '  -1.60%  '
